# Update the existing "helix jump" app row to the new "strechy" app/keyword,
# then append a brand-new review row for that same app.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 2: appid/keyword text fix (typo'd app name + renamed keyword)
$ws.Range("A2").Value = "com.singelton.strechy"
$ws.Range("B2").Value = "plank game"

# 2) New row 4: another review for the same app
$ws.Range("A4").Value = "com.singelton.strechy"
$ws.Range("B4").Value = "plank game"
$ws.Range("C4").Value = "armonravid2@gmail.com"
$ws.Range("D4").Value = "armonravid@gmail.com"
$ws.Range("E4").Value = "27/5/2019 15:59"
$ws.Range("F4").Value = "Great Game! I love it. Level 30 is so difficult!!"

# Recovery/email columns are mailto hyperlinks, same as the rows above them.
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:armonravid2@gmail.com", "", "", "armonravid2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:armonravid@gmail.com", "", "", "armonravid@gmail.com")

# Carry the same look-and-feel as the row above into the new row (Hyperlinks.Add
# otherwise stamps its own blue/underline style onto C4/D4). F4 (the review
# text) keeps the plain default look, same as every other unstyled column.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the selection left behind after typing the new row (one below it).
$ws.Range("F5").Select()
